# Edit script for Vocabulary/vocabu.xlsx
# Adds new vocabulary words (word/synonyms/meaning in Bengali) to several
# existing rows, two brand-new rows (59, 60), applies wrap-text styling to
# the four rows that contain longer two-line Bengali entries, widens
# columns C and D to fit the new content, and updates the sheet's view
# (selection / scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new cell values (in original authoring order, so that the
#     shared-strings table indices line up with the source workbook) ---
$ws.Range("A59").Value = "hypotheses"
$ws.Range("C59").Value = "guess, supposition, inference"
$ws.Range("C44").Value = "colleague,partner, cooperator"
$ws.Range("D43").Value = "have in common with"
$ws.Range("C43").Value = "likeness, take after"
$ws.Range("D24").Value = "ভক্তি`n"
$ws.Range("C24").Value = "piety, worship, adoration"
$ws.Range("C25").Value = "prominent,distinguised, noted"
$ws.Range("D27").Value = "মস্তকবিশিষ্ট`n"
$ws.Range("D28").Value = "সিনকোনা গাছের ছাল থেকে প্রাপ্ত উপক্ষার"
$ws.Range("D39").Value = "দিগন্ত"
$ws.Range("C39").Value = "skyline, azimuth"
$ws.Range("C41").Value = "erase,efface,mob,swob"
$ws.Range("D42").Value = "আশাবাদ"
$ws.Range("D21").Value = "হুমড়ি"
$ws.Range("C21").Value = "peck, nibbling, pecking, reprimand"
$ws.Range("D57").Value = "প্রত্যাশা"
$ws.Range("C57").Value = "hope, expectation, prospect"
$ws.Range("D56").Value = "খসড়া`n"
$ws.Range("C56").Value = "checkers, draughts"
$ws.Range("D55").Value = "উচ্চতর"
$ws.Range("C55").Value = "upper,best,excellent,beneficial"
$ws.Range("D47").Value = "লক্ষণীয়"
$ws.Range("C47").Value = "noticeable, remarkable"
$ws.Range("D50").Value = "জ্যোতির্বিদ্যা-সংক্রান্ত`n"
$ws.Range("A60").Value = "scepticism"
$ws.Range("D60").Value = "সংশয়বাদ"
$ws.Range("D58").Value = "অনুরূপ"
$ws.Range("C58").Value = "likeness"
$ws.Range("D52").Value = "ক্ষয়িত"

# --- Apply "wrap text" formatting + taller row height to the four cells
#     that hold two-line Bengali definitions ---
$wrapCells = @("D24", "D27", "D56", "D50")
foreach ($addr in $wrapCells) {
    $ws.Range($addr).WrapText = $true
}
$wrapRows = @(24, 27, 50, 56)
foreach ($r in $wrapRows) {
    $ws.Rows.Item($r).RowHeight = 30
}

# --- Widen columns C and D to fit the new synonym / meaning text ---
$ws.Columns.Item(3).ColumnWidth = 31
$ws.Columns.Item(4).ColumnWidth = 40.5

# --- Update the sheet view: scroll down to row 40 and select D52, matching
#     where the author ended up working in the sheet ---
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("D52").Select()
